$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 4 (pushes existing rows 4..42 down to 5..43)
$ws.Rows.Item(4).Insert()

# Fill the newly inserted row 4 with the new weekly record.
# Constant columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical across all
# data rows in this sheet, so copy them from row 5 (the row that used to
# be row 4 before the insert).
$ws.Range("A4").Value = $ws.Range("A5").Value()
$ws.Range("B4").Value = $ws.Range("B5").Value()
$ws.Range("C4").Value = $ws.Range("C5").Value()
$ws.Range("D4").Value = 44490
$ws.Range("E4").Value = $ws.Range("E5").Value()
$ws.Range("F4").Value = $ws.Range("F5").Value()
$ws.Range("G4").Value = $ws.Range("G5").Value()
$ws.Range("H4").Value = $ws.Range("H5").Value()
$ws.Range("I4").Value = $ws.Range("I5").Value()
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17000
$ws.Range("N4").Value = $ws.Range("N5").Value()
$ws.Range("O4").Value = $ws.Range("O5").Value()
$ws.Range("P4").Value = 680
$ws.Range("Q4").Value = $ws.Range("Q5").Value()
$ws.Range("R4").Value = $ws.Range("R5").Value()
